$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 43: FraxShare -> Frax
$ws.Range("B43").Value = "Frax"
$ws.Range("C43").Value = "https://coinranking.com/coin/KfWtaeV1W+frax-frax"

# Row 44: Frax -> FraxShare
$ws.Range("B44").Value = "FraxShare"
$ws.Range("C44").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"

# Price (D) and Volume(1h) (E) updates — use Text format to avoid Excel auto-converting
# the numeric-looking strings into actual numbers (the source data stores these as text).
$dCells = @(
    @{Cell="D2"; Value="27.888.70"},
    @{Cell="D3"; Value="1.745.55"},
    @{Cell="D5"; Value="334.23"},
    @{Cell="D7"; Value="0.3834"},
    @{Cell="D8"; Value="0.3376"},
    @{Cell="D9"; Value="45.85"},
    @{Cell="D10"; Value="1.108"},
    @{Cell="D11"; Value="0.07162"},
    @{Cell="D12"; Value="1.002"},
    @{Cell="D13"; Value="22.31"},
    @{Cell="D14"; Value="6.122"},
    @{Cell="D15"; Value="1.752.04"},
    @{Cell="D16"; Value="7.067"},
    @{Cell="D17"; Value="0.00001052"},
    @{Cell="D18"; Value="0.06596"},
    @{Cell="D19"; Value="78.81"},
    @{Cell="D21"; Value="16.65"},
    @{Cell="D22"; Value="6.158"},
    @{Cell="D23"; Value="27.956.70"},
    @{Cell="D24"; Value="11.58"},
    @{Cell="D26"; Value="153.25"},
    @{Cell="D27"; Value="19.78"},
    @{Cell="D28"; Value="2.287"},
    @{Cell="D29"; Value="1.950.15"},
    @{Cell="D30"; Value="1.277"},
    @{Cell="D31"; Value="130.35"},
    @{Cell="D32"; Value="4.023"},
    @{Cell="D33"; Value="5.772"},
    @{Cell="D34"; Value="0.08769"},
    @{Cell="D35"; Value="12.08"},
    @{Cell="D37"; Value="0.6529"},
    @{Cell="D38"; Value="0.02269"},
    @{Cell="D39"; Value="5.099"},
    @{Cell="D40"; Value="0.06077"},
    @{Cell="D41"; Value="0.2077"},
    @{Cell="D42"; Value="1.202"},
    @{Cell="D43"; Value="1.000"},
    @{Cell="D44"; Value="7.907"},
    @{Cell="D45"; Value="13.61"},
    @{Cell="D47"; Value="0.5979"},
    @{Cell="D48"; Value="126.41"},
    @{Cell="D49"; Value="1.990"},
    @{Cell="D50"; Value="1.161"},
    @{Cell="D51"; Value="1.104"}
)

foreach ($item in $dCells) {
    $rng = $ws.Range($item.Cell)
    $rng.NumberFormat = "@"
    $rng.Value = $item.Value
}

$eCells = @(
    @{Cell="E2"; Value="  +1.36%  "},
    @{Cell="E3"; Value="  -0.97%  "},
    @{Cell="E4"; Value="  +0.02%  "},
    @{Cell="E5"; Value="  -0.49%  "},
    @{Cell="E6"; Value="  -0.13%  "},
    @{Cell="E7"; Value="  +0.04%  "},
    @{Cell="E8"; Value="  -0.79%  "},
    @{Cell="E9"; Value="  -2.23%  "},
    @{Cell="E10"; Value="  -2.51%  "},
    @{Cell="E11"; Value="  -2.97%  "},
    @{Cell="E12"; Value="  -0.04%  "},
    @{Cell="E13"; Value="  -0.55%  "},
    @{Cell="E14"; Value="  -3.37%  "},
    @{Cell="E15"; Value="  -0.68%  "},
    @{Cell="E16"; Value="  +0.43%  "},
    @{Cell="E17"; Value="  -1.73%  "},
    @{Cell="E18"; Value="  -1.03%  "},
    @{Cell="E19"; Value="  -3.84%  "},
    @{Cell="E20"; Value="  -0.04%  "},
    @{Cell="E22"; Value="  -3.23%  "},
    @{Cell="E23"; Value="  +1.55%  "},
    @{Cell="E24"; Value="  -3.33%  "},
    @{Cell="E25"; Value="  +0.39%  "},
    @{Cell="E26"; Value="  +0.70%  "},
    @{Cell="E27"; Value="  -3.75%  "},
    @{Cell="E28"; Value="  -5.07%  "},
    @{Cell="E29"; Value="  -0.70%  "},
    @{Cell="E30"; Value="  -10.21%  "},
    @{Cell="E31"; Value="  -3.00%  "},
    @{Cell="E32"; Value="  +1.64%  "},
    @{Cell="E33"; Value="  -5.04%  "},
    @{Cell="E34"; Value="  -0.14%  "},
    @{Cell="E35"; Value="  -4.70%  "},
    @{Cell="E36"; Value="  +1.93%  "},
    @{Cell="E37"; Value="  -3.46%  "},
    @{Cell="E38"; Value="  -5.67%  "},
    @{Cell="E39"; Value="  -3.81%  "},
    @{Cell="E40"; Value="  -3.48%  "},
    @{Cell="E41"; Value="  -4.40%  "},
    @{Cell="E43"; Value="  -0.07%  "},
    @{Cell="E44"; Value="  -3.85%  "},
    @{Cell="E45"; Value="  -3.45%  "},
    @{Cell="E46"; Value="  +0.29%  "},
    @{Cell="E47"; Value="  -3.93%  "},
    @{Cell="E48"; Value="  -3.01%  "},
    @{Cell="E49"; Value="  -3.81%  "},
    @{Cell="E50"; Value="  +1.44%  "},
    @{Cell="E51"; Value="  +4.48%  "}
)

foreach ($item in $eCells) {
    $ws.Range($item.Cell).Value = $item.Value
}
